# Gantt chart update — Project Plan & Gantt Chart Update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Planner")

# --- Update Percent Complete / Actual Start / Actual Duration values ---
# Row 17 (Activity 14): mark fully complete
$ws.Range("G17").Value = 1

# Row 18 (Activity 15): mark fully complete
$ws.Range("G18").Value = 1

# Row 19 (Activity 16): mark fully complete
$ws.Range("G19").Value = 1

# Row 20 (Activity 17): Actual Start 8, Actual Duration 8, 100% complete
$ws.Range("E20").Value = 8
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 1

# Row 21 (Activity 18): Actual Start 10, Actual Duration 16, 100% complete
$ws.Range("E21").Value = 10
$ws.Range("F21").Value = 16
$ws.Range("G21").Value = 1

# Row 22 (Activity 19): Actual Start 18, Actual Duration 15, 100% complete
$ws.Range("E22").Value = 18
$ws.Range("F22").Value = 15
$ws.Range("G22").Value = 1

# Row 23 (Activity 20): Actual Start 16, Actual Duration 18, 100% complete
$ws.Range("E23").Value = 16
$ws.Range("F23").Value = 18
$ws.Range("G23").Value = 1

# Row 24 (Activity 21): Actual Start 17, Actual Duration 17, 100% complete
$ws.Range("E24").Value = 17
$ws.Range("F24").Value = 17
$ws.Range("G24").Value = 1

# Row 25 (Activity 22): Actual Start 18, Actual Duration 17, 100% complete
$ws.Range("E25").Value = 18
$ws.Range("F25").Value = 17
$ws.Range("G25").Value = 1

# Row 26 (Activity 23): Actual Start 20, Actual Duration 18, 100% complete
$ws.Range("E26").Value = 20
$ws.Range("F26").Value = 18
$ws.Range("G26").Value = 1

# Row 27 (Activity 24): Actual Start 20, Actual Duration 20, 100% complete
$ws.Range("E27").Value = 20
$ws.Range("F27").Value = 20
$ws.Range("G27").Value = 1

# Row 28 (Activity 25): Actual Start 21, Actual Duration 20, 100% complete
$ws.Range("E28").Value = 21
$ws.Range("F28").Value = 20
$ws.Range("G28").Value = 1

# Row 29 (Activity 26): Actual Start 22, Actual Duration 20, 100% complete
$ws.Range("E29").Value = 22
$ws.Range("F29").Value = 20
$ws.Range("G29").Value = 1

# Row 30 (Activity 27): Actual Start 21, Actual Duration 25, 100% complete
$ws.Range("E30").Value = 21
$ws.Range("F30").Value = 25
$ws.Range("G30").Value = 1

# Row 31 (Activity 28): Actual Start 21, Actual Duration 25, 100% complete
$ws.Range("E31").Value = 21
$ws.Range("F31").Value = 25
$ws.Range("G31").Value = 1

# Row 32 (Activity 29): Actual Start 22, Actual Duration 27, 100% complete
$ws.Range("E32").Value = 22
$ws.Range("F32").Value = 27
$ws.Range("G32").Value = 1

# Row 33 (Activity 30): Actual Start 24, Actual Duration 26, 100% complete
$ws.Range("E33").Value = 24
$ws.Range("F33").Value = 26
$ws.Range("G33").Value = 1

# --- View state: scroll down, adjust zoom, reselect active cell ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 75
$ws.Range("G33").Select()
